$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "29.479.52"
Set-TextValue "E2" "  +0.77%  "

Set-TextValue "D3" "1.971.53"
Set-TextValue "E3" "  +3.82%  "

Set-TextValue "E4" "  +0.29%  "

Set-TextValue "D5" "326.67"
Set-TextValue "E5" "  -0.19%  "

Set-TextValue "D6" "1.003"
Set-TextValue "E6" "  +0.17%  "

Set-TextValue "D7" "0.4658"
Set-TextValue "E7" "  +0.56%  "

Set-TextValue "D8" "0.3912"
Set-TextValue "E8" "  -0.32%  "

Set-TextValue "D9" "46.15"

Set-TextValue "D10" "0.07933"
Set-TextValue "E10" "  +0.70%  "

Set-TextValue "D11" "0.9891"
Set-TextValue "E11" "  +0.08%  "

Set-TextValue "D12" "22.72"
Set-TextValue "E12" "  +4.25%  "

Set-TextValue "D13" "1.976.61"
Set-TextValue "E13" "  +2.96%  "

Set-TextValue "D14" "7.176"
Set-TextValue "E14" "  +1.54%  "

Set-TextValue "D15" "5.829"
Set-TextValue "E15" "  +1.69%  "

Set-TextValue "D16" "0.07064"
Set-TextValue "E16" "  +0.98%  "

Set-TextValue "D17" "87.64"
Set-TextValue "E17" "  -0.80%  "

Set-TextValue "E18" "  +0.29%  "

Set-TextValue "D19" "0.000009936"
Set-TextValue "E19" "  -0.26%  "

Set-TextValue "D20" "17.26"
Set-TextValue "E20" "  +1.25%  "

Set-TextValue "D21" "1.003"
Set-TextValue "E21" "  +0.09%  "

Set-TextValue "D22" "29.488.67"
Set-TextValue "E22" "  +0.75%  "

Set-TextValue "B23" "Uniswap"
Set-TextValue "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D23" "5.533"
Set-TextValue "E23" "  +4.34%  "

Set-TextValue "B24" "Cosmos"
Set-TextValue "C24" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D24" "11.14"
Set-TextValue "E24" "  +0.49%  "

Set-TextValue "B25" "WrappedliquidstakedEther2.0"
Set-TextValue "C25" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D25" "2.224.23"
Set-TextValue "E25" "  +3.68%  "

Set-TextValue "B26" "Toncoin"
Set-TextValue "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "2.108"
Set-TextValue "E26" "  +0.39%  "

Set-TextValue "B27" "Monero"
Set-TextValue "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "158.50"
Set-TextValue "E27" "  +1.65%  "

Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "19.52"
Set-TextValue "E28" "  +0.46%  "

Set-TextValue "B29" "InternetComputer(DFINITY)"
Set-TextValue "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D29" "5.765"
Set-TextValue "E29" "  -4.46%  "

Set-TextValue "B30" "BitcoinCash"
Set-TextValue "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D30" "119.44"
Set-TextValue "E30" "  +0.87%  "

Set-TextValue "B31" "LidoDAOToken"
Set-TextValue "C31" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D31" "1.903"
Set-TextValue "E31" "  +0.74%  "

Set-TextValue "B32" "Stellar"
Set-TextValue "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D32" "0.09412"
Set-TextValue "E32" "  +0.61%  "

Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.8901"
Set-TextValue "E33" "  -1.53%  "

Set-TextValue "B34" "Filecoin"
Set-TextValue "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "5.230"
Set-TextValue "E34" "  -0.36%  "

Set-TextValue "B35" "ARBITRUM"
Set-TextValue "C35" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "1.319"
Set-TextValue "E35" "  -0.39%  "

Set-TextValue "B36" "HuobiToken"
Set-TextValue "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "3.168"
Set-TextValue "E36" "  -1.41%  "

Set-TextValue "B37" "Hedera"
Set-TextValue "C37" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D37" "0.05818"
Set-TextValue "E37" "  +0.75%  "

Set-TextValue "B38" "TrustWalletToken"
Set-TextValue "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.171"
Set-TextValue "E38" "  -1.26%  "

Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.02100"
Set-TextValue "E39" "  +0.61%  "

Set-TextValue "B40" "FraxShare"
Set-TextValue "C40" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "7.746"
Set-TextValue "E40" "  +0.29%  "

Set-TextValue "B41" "TheSandbox"
Set-TextValue "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.5722"
Set-TextValue "E41" "  +0.35%  "

Set-TextValue "B42" "PEPE"
Set-TextValue "C42" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D42" "0.000003088"
Set-TextValue "E42" "  +47.43%  "

Set-TextValue "B43" "Algorand"
Set-TextValue "C43" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D43" "0.1795"
Set-TextValue "E43" "  +0.54%  "

Set-TextValue "B44" "Aptos"
Set-TextValue "C44" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D44" "9.637"
Set-TextValue "E44" "  -0.70%  "

Set-TextValue "B45" "MXToken"
Set-TextValue "C45" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D45" "2.757"
Set-TextValue "E45" "  +7.13%  "

Set-TextValue "B46" "Decentraland"
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5340"
Set-TextValue "E46" "  -0.28%  "

Set-TextValue "D47" "11.70"
Set-TextValue "E47" "  -1.46%  "

Set-TextValue "D48" "2.195"
Set-TextValue "E48" "  +0.69%  "

Set-TextValue "B49" "Cronos"
Set-TextValue "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.06920"
Set-TextValue "E49" "  -1.51%  "

Set-TextValue "B50" "NEARProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "1.828"
Set-TextValue "E50" "  -1.19%  "

Set-TextValue "B51" "Quant"
Set-TextValue "C51" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D51" "113.59"
Set-TextValue "E51" "  +0.52%  "
